$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the "X" value out of E5 (leave the cell's formatting/style intact)
$ws.Range("E5").ClearContents()

# Reflect the user's final selection on the sheet
$ws.Range("F22").Select()
